# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to the newly scraped values.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $ws.Range("F2").Value  = 1250
        $ws.Range("F3").Value  = 17102
        $ws.Range("F4").Value  = 51
        $ws.Range("F5").Value  = 1666
        $ws.Range("F8").Value  = 1035
        $ws.Range("F10").Value = 240
        $ws.Range("F12").Value = 11859
        $ws.Range("F13").Value = 34
        $ws.Range("F14").Value = 44
        $ws.Range("F15").Value = 11532
        $ws.Range("F16").Value = 4706
        $ws.Range("F17").Value = 497
        $ws.Range("F18").Value = 55
        $ws.Range("F24").Value = 44
    }
    elseif ($sheetName -eq "全部类型") {
        $ws.Range("F2").Value  = 1250
        $ws.Range("F3").Value  = 17102
        $ws.Range("F4").Value  = 51
        $ws.Range("F5").Value  = 1666
        $ws.Range("F8").Value  = 1035
        $ws.Range("F10").Value = 240
        $ws.Range("F14").Value = 11859
        $ws.Range("F15").Value = 34
        $ws.Range("F16").Value = 44
        $ws.Range("F17").Value = 11533
        $ws.Range("F18").Value = 4706
        $ws.Range("F19").Value = 497
        $ws.Range("F20").Value = 55
        $ws.Range("F26").Value = 44
    }
}
